# BigChartProjeto.xlsx - "atualizando big chart conforme recomendado"
#
# A new measurement row (22/09/2010 -> serial 40429) was inserted at the top
# of the Plan1 data table (all metrics start at 0 on that date), pushing the
# existing rows down by one, and a new tracked metric "Controladoras" was
# added as column G with its own header, data and a matching line in the
# chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# ---------------------------------------------------------------------
# 1) Worksheet data - write the final table state directly (header row +
#    8 date rows x 7 columns). This reproduces the effect of inserting a
#    new first data row and appending a new "Controladoras" column.
# ---------------------------------------------------------------------

# Header row
$ws.Cells.Item(1,1).Value = "Data"
$ws.Cells.Item(1,2).Value = "Classes de modelo"
$ws.Cells.Item(1,3).Value = "Testes de Unidade"
$ws.Cells.Item(1,4).Value = "Testes de Aceitação"
$ws.Cells.Item(1,5).Value = "Páginas GSP"
$ws.Cells.Item(1,6).Value = "User Stories"
$ws.Cells.Item(1,7).Value = "Controladoras"

# New first measurement row (all metrics at zero on the earliest date)
$ws.Cells.Item(2,1).Value = 40429
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,7).Font.Bold = $false
$ws.Cells.Item(2,7).HorizontalAlignment = -4108

# Former row 2 (now row 3) keeps its values and gains "Controladoras"=3
$ws.Cells.Item(3,1).Value = 40443
$ws.Cells.Item(3,2).Value = 2
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = 5
$ws.Cells.Item(3,5).Value = 15
$ws.Cells.Item(3,6).Value = 5
$ws.Cells.Item(3,7).Value = 3

# Remaining dates shift down by one row, keeping their (empty) metric cells
$ws.Cells.Item(4,1).Value = 40457
$ws.Cells.Item(4,2).Value = $null
$ws.Cells.Item(4,3).Value = $null
$ws.Cells.Item(4,4).Value = $null
$ws.Cells.Item(4,5).Value = $null
$ws.Cells.Item(4,6).Value = $null

$ws.Cells.Item(5,1).Value = 40471
$ws.Cells.Item(5,2).Value = $null
$ws.Cells.Item(5,3).Value = $null
$ws.Cells.Item(5,4).Value = $null
$ws.Cells.Item(5,5).Value = $null
$ws.Cells.Item(5,6).Value = $null

$ws.Cells.Item(6,1).Value = 40485
$ws.Cells.Item(6,2).Value = $null
$ws.Cells.Item(6,3).Value = $null
$ws.Cells.Item(6,4).Value = $null
$ws.Cells.Item(6,5).Value = $null
$ws.Cells.Item(6,6).Value = $null

$ws.Cells.Item(7,1).Value = 40499
$ws.Cells.Item(7,2).Value = $null
$ws.Cells.Item(7,3).Value = $null
$ws.Cells.Item(7,4).Value = $null
$ws.Cells.Item(7,5).Value = $null
$ws.Cells.Item(7,6).Value = $null

$ws.Cells.Item(8,1).Value = 40513
$ws.Cells.Item(8,2).Value = $null
$ws.Cells.Item(8,3).Value = $null
$ws.Cells.Item(8,4).Value = $null
$ws.Cells.Item(8,5).Value = $null
$ws.Cells.Item(8,6).Value = $null

# Column G width, similar to the other bestFit columns
$ws.Columns.Item(7).ColumnWidth = 18.36

# Selection moved by the editor onto the new column while working on it
$ws.Range("G7").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Chart - add the "Controladoras" series plotted against the same
#    date axis, sourced from the new column G.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$newSeries = $chart.SeriesCollection().NewSeries()
$newSeries.Name = "Controladoras"
$newSeries.Values = "=Plan1!`$G`$2:`$G`$8"
$newSeries.XValues = "=Plan1!`$A`$2:`$A`$8"
